$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "New storms headed for Midwest, Northeast"
$ws.Range("B2").Value = "2007-12-15T03:21:08UTC"
$ws.Range("C2").Value = 348
$ws.Range("D2").Value = "day_31_beyond"
$ws.Range("E2").Value = "http://www.nbcnews.com/id/22239846"

$ws.Range("A3").Value = "Preliminary Information for the December 2007 Ice Storm in Oklahoma"
$ws.Range("B3").Value = "2007-12-08T00:00:00UTC"
$ws.Range("C3").Value = 341
$ws.Range("D3").Value = "day_31_beyond"
$ws.Range("E3").Value = "http://www.srh.noaa.gov/oun/wxevents/20071208/photos.php"

$ws.Range("A4").Value = "Ice coats nation’s midsection"
$ws.Range("B4").Value = "2007-12-11T03:56:00UTC"
$ws.Range("C4").Value = 344
$ws.Range("D4").Value = "day_31_beyond"
$ws.Range("E4").Value = "http://www.nbcnews.com/id/22173398"

$ws.Range("A5").Value = "Northeast could see double winter whammy"
$ws.Range("B5").Value = "2007-12-13T01:55:23UTC"
$ws.Range("C5").Value = 346
$ws.Range("D5").Value = "day_31_beyond"
$ws.Range("E5").Value = "http://www.nbcnews.com/id/22197793"

$ws.Range("A6").Value = "Ice Storm Brings Misery, Death to Midwest"
$ws.Range("B6").Value = "2009-02-26T00:00:00UTC"
$ws.Range("C6").Value = 787
$ws.Range("D6").Value = "day_31_beyond"
$ws.Range("E6").Value = "https://abcnews.go.com/US/wireStory?id=3983346"

$ws.Range("A7").Value = "Blizzard that walloped Ont., Que., heads to Eastern Canada"
$ws.Range("B7").Value = "2007-12-16T00:00:00UTC"
$ws.Range("C7").Value = 349
$ws.Range("D7").Value = "day_31_beyond"
$ws.Range("E7").Value = "https://www.cbc.ca/news/canada/blizzard-that-walloped-ont-que-heads-to-eastern-canada-1.630843"

$ws.Range("A8").Value = "Ice and Snow Storm of December 13th"
$ws.Range("B8").Value = "1-01-01T00:00:00UTC"
$ws.Range("C8").Value = "unknown"
$ws.Range("D8").Value = "unknown"
$ws.Range("E8").Value = "http://www.erh.noaa.gov/ctp/features/2007/12_13/index.php"

$ws.Range("A9").Value = "March 23, 2007 Severe Weather Event"
$ws.Range("B9").Value = "1-01-01T00:00:00UTC"
$ws.Range("C9").Value = "unknown"
$ws.Range("D9").Value = "unknown"
$ws.Range("E9").Value = "http://www.srh.noaa.gov/abq/quickfeatures/Dec2007/RainSnowEventDec7-12.php"

$ws.Range("A10").Value = "December 16-17, 2007 Winter Storms"
$ws.Range("B10").Value = "1-01-01T00:00:00UTC"
$ws.Range("C10").Value = "unknown"
$ws.Range("D10").Value = "unknown"
$ws.Range("E10").Value = "http://www.erh.noaa.gov/btv/events/16-17Dec2007/16-17Dec2007.shtml"

$ws.Range("A11").Value = "NWS SGF Winter Summary"
$ws.Range("B11").Value = "1-01-01T00:00:00UTC"
$ws.Range("C11").Value = "unknown"
$ws.Range("D11").Value = "unknown"
$ws.Range("E11").Value = "http://www.crh.noaa.gov/sgf/?n=icestormsummarydec07"

$ws.Range("A12").Value = "National Weather Service Forecast Office"
$ws.Range("B12").Value = "1-01-01T00:00:00UTC"
$ws.Range("C12").Value = "unknown"
$ws.Range("D12").Value = "unknown"
$ws.Range("E12").Value = "https://web.archive.org/web/20080511204226/http://www.erh.noaa.gov/btv/events/IceStorm1998/ice98.shtml"

$ws.Range("A13").Value = "Oklahomans catch a break"
$ws.Range("B13").Value = "1-01-01T00:00:00UTC"
$ws.Range("C13").Value = "unknown"
$ws.Range("D13").Value = "unknown"
$ws.Range("E13").Value = "https://web.archive.org/web/20071215171743/http://edition.cnn.com/2007/US/weather/12/15/winter.storm.ap/index.html"

$ws.Range("A14").Value = "New rain, ice hamper storm recovery effort"
$ws.Range("B14").Value = "1-01-01T00:00:00UTC"
$ws.Range("C14").Value = "unknown"
$ws.Range("D14").Value = "unknown"
$ws.Range("E14").Value = "https://web.archive.org/web/20071213204855/http://edition.cnn.com/2007/US/weather/12/12/winter.storm.ap/index.html"

$ws.Range("A15").Value = "NOAA's National Weather Service"
$ws.Range("B15").Value = "1-01-01T00:00:00UTC"
$ws.Range("C15").Value = "unknown"
$ws.Range("D15").Value = "unknown"
$ws.Range("E15").Value = "http://www.erh.noaa.gov/er/bgm/WeatherEvents/Snow/april162007/april162007.shtml"

$ws.Range("A16").Value = "NWS Des Moines, IA Winter Weather Information"
$ws.Range("B16").Value = "1-01-01T00:00:00UTC"
$ws.Range("C16").Value = "unknown"
$ws.Range("D16").Value = "unknown"
$ws.Range("E16").Value = "http://www.crh.noaa.gov/dmx/winterwx-post.php"

$ws.Range("A17").Value = "'Nasty' winter storms cause traffic deaths, delays"
$ws.Range("B17").Value = "1-01-01T00:00:00UTC"
$ws.Range("C17").Value = "unknown"
$ws.Range("D17").Value = "unknown"
$ws.Range("E17").Value = "https://web.archive.org/web/20071218174133/http://edition.cnn.com/2007/US/weather/12/17/winter.storm.ap/index.html"

$ws.Range("A18").Value = "Furious snow storm blows north, blankets Great Lakes states"
$ws.Range("B18").Value = "1-01-01T00:00:00UTC"
$ws.Range("C18").Value = "unknown"
$ws.Range("D18").Value = "unknown"
$ws.Range("E18").Value = "https://web.archive.org/web/20071217224528/http://edition.cnn.com/2007/US/weather/12/16/winter.storm.ap/index.html"

$ws.Range("A19").Value = "NOAA's NWS Forecast Office-Caribou, Maine-Snowfall Totals for the 11-12 Dec 2007 Winter Storm"
$ws.Range("B19").Value = "1-01-01T00:00:00UTC"
$ws.Range("C19").Value = "unknown"
$ws.Range("D19").Value = "unknown"
$ws.Range("E19").Value = "http://www.erh.noaa.gov/car/News_Items/2007-12-12_item001.htm"

$ws.Range("A20").Value = "Winter Storm Summary"
$ws.Range("B20").Value = "1-01-01T00:00:00UTC"
$ws.Range("C20").Value = "unknown"
$ws.Range("D20").Value = "unknown"
$ws.Range("E20").Value = "http://www.crh.noaa.gov/sgf/?n=snowsummary121507"

$ws.Range("A21").Value = "Midwest buried under heap of snow"
$ws.Range("B21").Value = "1-01-01T00:00:00UTC"
$ws.Range("C21").Value = "unknown"
$ws.Range("D21").Value = "unknown"
$ws.Range("E21").Value = "http://edition.cnn.com/2007/US/weather/12/15/winter.storm.ap/index.html"

$ws.Range("A22").Value = "Storm Prediction Center 20071215's Storm Reports"
$ws.Range("B22").Value = "1-01-01T00:00:00UTC"
$ws.Range("C22").Value = "unknown"
$ws.Range("D22").Value = "unknown"
$ws.Range("E22").Value = "http://www.spc.noaa.gov/climo/reports/071215_rpts.html"

$ws.Range("A23").Value = "From the Wichita National Weather Service Press"
$ws.Range("B23").Value = "1-01-01T00:00:00UTC"
$ws.Range("C23").Value = "unknown"
$ws.Range("D23").Value = "unknown"
$ws.Range("E23").Value = "http://www.crh.noaa.gov/ict/scripts/viewstory.php?STORY_NUMBER=2007121516"

$ws.Range("A24").Value = "Tampa Bay Area Aviation Weather Page"
$ws.Range("B24").Value = "1-01-01T00:00:00UTC"
$ws.Range("C24").Value = "unknown"
$ws.Range("D24").Value = "unknown"
$ws.Range("E24").Value = "http://www.srh.noaa.gov/tbw/html/tbw/getprodversionnew.php?pil=PNS&sid=TBW&max=10&"

$ws.Range("A25").Value = "Winter Weather Event"
$ws.Range("B25").Value = "1-01-01T00:00:00UTC"
$ws.Range("C25").Value = "unknown"
$ws.Range("D25").Value = "unknown"
$ws.Range("E25").Value = "http://www.ok.gov/OEM/Emergencies_&_Disasters/2007/Winter_Weather_Event_20071209_-_Master/"

$ws.Range("A26").Value = "snowfallTotals"
$ws.Range("B26").Value = "1-01-01T00:00:00UTC"
$ws.Range("C26").Value = "unknown"
$ws.Range("D26").Value = "unknown"
$ws.Range("E26").Value = "http://www.crh.noaa.gov/bou/?n=snowfalltotals_121107"

$ws.Range("A27").Value = "2007Dec10"
$ws.Range("B27").Value = "1-01-01T00:00:00UTC"
$ws.Range("C27").Value = "unknown"
$ws.Range("D27").Value = "unknown"
$ws.Range("E27").Value = "http://www.crh.noaa.gov/eax/?n=december10-11,2007icestorm"

$ws.Range("A28").Value = "NOAA's NWS Forecast Office-Caribou, Maine-Preliminary Totals for the Storm of 13-14 Dec. 2007"
$ws.Range("B28").Value = "1-01-01T00:00:00UTC"
$ws.Range("C28").Value = "unknown"
$ws.Range("D28").Value = "unknown"
$ws.Range("E28").Value = "http://www.erh.noaa.gov/car/News_Items/2007-12-14_item001.htm"

$ws.Range("A29").Value = "Dec 9-10 Ice Storm Summary"
$ws.Range("B29").Value = "1-01-01T00:00:00UTC"
$ws.Range("C29").Value = "unknown"
$ws.Range("D29").Value = "unknown"
$ws.Range("E29").Value = "http://www.srh.noaa.gov/tsa/weather-events/IceStormDec10/Index.htm"
